# Forecast sheet edits + FIM run
# Re-applies the updated forward-indexed-model (FIM) projection values
# for rows 219-260 on the single worksheet in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 219
$ws.Range("AD219").Value = 4417.2

# Row 220
$ws.Range("C220").Value = 29354.3
$ws.Range("E220").Value = 919.4
$ws.Range("G220").Value = 4424.9
$ws.Range("H220").Value = 1893.4
$ws.Range("I220").Value = 3139.2
$ws.Range("J220").Value = 483.1
$ws.Range("K220").Value = 2152.299
$ws.Range("L220").Value = 164.8
$ws.Range("M220").Value = 221.6
$ws.Range("T220").Value = 0.0037233826809977
$ws.Range("U220").Value = 0.00608116293155669
$ws.Range("V220").Value = 0.00575477607265906
$ws.Range("W220").Value = 0.00523456712027648
$ws.Range("X220").Value = 0.00792756208803924
$ws.Range("AD220").Value = 4465.4
$ws.Range("AE220").Value = 2483.3
$ws.Range("AL220").Value = 278.163

# Row 221
$ws.Range("C221").Value = 29645.7262081059
$ws.Range("E221").Value = 927.689732874973
$ws.Range("H221").Value = 1885.09320921534
$ws.Range("I221").Value = 3167.88900514805
$ws.Range("J221").Value = 492.801912246684
$ws.Range("K221").Value = 2160.527
$ws.Range("L221").Value = 163.15466849793
$ws.Range("N221").Value = 647.018687447845
$ws.Range("T221").Value = 0.00545171502706143
$ws.Range("AD221").Value = 4505.15020255048
$ws.Range("AE221").Value = 2522.33520327819
$ws.Range("AK221").Value = 1763.0346223627
$ws.Range("AL221").Value = 280.671045427128

# Row 222
$ws.Range("C222").Value = 29949.8187152429
$ws.Range("E222").Value = 928.991094574193
$ws.Range("H222").Value = 1907.3318350055
$ws.Range("I222").Value = 3195.56015057414
$ws.Range("J222").Value = 502.698664280664
$ws.Range("K222").Value = 2202.093889974
$ws.Range("L222").Value = 164.648939314487
$ws.Range("N222").Value = 647.926324245673
$ws.Range("AD222").Value = 4578.04931402487
$ws.Range("AE222").Value = 2548.44772380158
$ws.Range("AK222").Value = 1778.08535540041
$ws.Range("AL222").Value = 281.06477032852

# Row 223
$ws.Range("C223").Value = 30261.308341014
$ws.Range("E223").Value = 932.294281821561
$ws.Range("H223").Value = 1930.4539814257
$ws.Range("I223").Value = 3223.36403051996
$ws.Range("J223").Value = 512.794168994753
$ws.Range("K223").Value = 2208.793889974
$ws.Range("L223").Value = 166.13773661157
$ws.Range("N223").Value = 650.230137472713
$ws.Range("AD223").Value = 4609.71515475437
$ws.Range("AE223").Value = 2574.48894008494
$ws.Range("AK223").Value = 1794.71149812467
$ws.Range("AL223").Value = 282.064144348848

# Row 224
$ws.Range("C224").Value = 30566.7181432502
$ws.Range("E224").Value = 935.602102774806
$ws.Range("H224").Value = 1954.158923269
$ws.Range("I224").Value = 3250.16069012155
$ws.Range("J224").Value = 523.092417862893
$ws.Range("K224").Value = 2215.493889974
$ws.Range("L224").Value = 167.08465548067
$ws.Range("N224").Value = 652.537182485326
$ws.Range("AD224").Value = 4641.62782503496
$ws.Range("AE224").Value = 2599.52236547798
$ws.Range("AK224").Value = 1811.522377287
$ws.Range("AL224").Value = 283.064920289481

# Row 225
$ws.Range("C225").Value = 30856.5230650801
$ws.Range("E225").Value = 941.843153399173
$ws.Range("H225").Value = 1968.11289188136
$ws.Range("I225").Value = 3276.93438119768
$ws.Range("J225").Value = 521.771537588453
$ws.Range("K225").Value = 2205.954889974
$ws.Range("L225").Value = 167.166758272789
$ws.Range("N225").Value = 656.890013221912
$ws.Range("AD225").Value = 4673.78951251855
$ws.Range("AE225").Value = 2624.24639104833
$ws.Range("AK225").Value = 1834.11720688667
$ws.Range("AL225").Value = 284.953140177262

# Row 226
$ws.Range("C226").Value = 31126.264569245
$ws.Range("E226").Value = 948.11249448117
$ws.Range("H226").Value = 1980.71518106142
$ws.Range("I226").Value = 3304.82349595586
$ws.Range("J226").Value = 520.453992718313
$ws.Range("K226").Value = 2249.7275234174
$ws.Range("L226").Value = 166.433306663197
$ws.Range("N226").Value = 661.262575183405
$ws.Range("AD226").Value = 4755.93289067343
$ws.Range("AE226").Value = 2646.38928336033
$ws.Range("AK226").Value = 1857.01889003517
$ws.Range("AL226").Value = 286.849919297766

# Row 227
$ws.Range("C227").Value = 31399.8566283152
$ws.Range("E227").Value = 954.410254260423
$ws.Range("H227").Value = 1995.67745272291
$ws.Range("I227").Value = 3332.32958624848
$ws.Range("J227").Value = 519.139774830117
$ws.Range("K227").Value = 2255.6995234174
$ws.Range("L227").Value = 166.312889234756
$ws.Range("N227").Value = 665.654957810736
$ws.Range("AD227").Value = 4788.84607274818
$ws.Range("AE227").Value = 2669.39124256
$ws.Range("AK227").Value = 1880.23203537275
$ws.Range("AL227").Value = 288.755296449687

# Row 228
$ws.Range("C228").Value = 31684.696360925
$ws.Range("E228").Value = 960.736561557861
$ws.Range("H228").Value = 2009.40022998616
$ws.Range("I228").Value = 3360.14586941614
$ws.Range("J228").Value = 517.828875522779
$ws.Range("K228").Value = 2263.3995234174
$ws.Range("L228").Value = 166.756244312196
$ws.Range("N228").Value = 670.06725095027
$ws.Range("AD228").Value = 4822.01619301344
$ws.Range("AE228").Value = 2692.90294087333
$ws.Range("AK228").Value = 1903.76132307562
$ws.Range("AL228").Value = 290.669310607591

# Row 229
$ws.Range("C229").Value = 31972.5760053022

# Row 230
$ws.Range("C230").Value = 32256.3011035973

# Row 231
$ws.Range("C231").Value = 32537.3916116939

# Row 232
$ws.Range("C232").Value = 32821.8260227347

# Row 233
$ws.Range("C233").Value = 33110.0096582887

# Row 234
$ws.Range("C234").Value = 33402.5505007094

# Row 235
$ws.Range("C235").Value = 33702.0831401952

# Row 236
$ws.Range("C236").Value = 34008.1009247848

# Row 237
$ws.Range("C237").Value = 34320.5025240862

# Row 238
$ws.Range("C238").Value = 34638.3759645689

# Row 239
$ws.Range("C239").Value = 34961.8225766253

# Row 240
$ws.Range("C240").Value = 35290.0317171173

# Row 241
$ws.Range("C241").Value = 35622.0914125148

# Row 242
$ws.Range("C242").Value = 35957.1910196797

# Row 243
$ws.Range("C243").Value = 36295.8371905733

# Row 244
$ws.Range("C244").Value = 36638.1312555879

# Row 245
$ws.Range("C245").Value = 36983.9718843311

# Row 246
$ws.Range("C246").Value = 37332.953755234

# Row 247
$ws.Range("C247").Value = 37685.1781986889

# Row 248
$ws.Range("C248").Value = 38040.8478754802

# Row 249
$ws.Range("C249").Value = 38398.9494816855

# Row 250
$ws.Range("C250").Value = 38759.4830173047

# Row 251
$ws.Range("C251").Value = 39123.0564646914

# Row 252
$ws.Range("C252").Value = 39489.87248463

# Row 253
$ws.Range("C253").Value = 39859.5257555516

# Row 254
$ws.Range("C254").Value = 40231.9149470639

# Row 255
$ws.Range("C255").Value = 40607.2427199514

# Row 256
$ws.Range("C256").Value = 40985.7117349985

# Row 257
$ws.Range("C257").Value = 41367.4233225977

# Row 258
$ws.Range("C258").Value = 41752.6814739255

# Row 259
$ws.Range("C259").Value = 42141.2835281975

# Row 260
$ws.Range("C260").Value = 42533.4321461982
